$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (style s="1") from G1 into H1:V1
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:V1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Set header labels for H1:V1
$ws.Range("H1").Value = "EoD_Bestand_noSiBe_log1p"
$ws.Range("I1").Value = "DeficitPos_log1p"
$ws.Range("J1").Value = "DemandMean_100_log1p"
$ws.Range("K1").Value = "DemandMax_100_log1p"
$ws.Range("L1").Value = "DemandMean_66_log1p"
$ws.Range("M1").Value = "DemandMax_66_log1p"
$ws.Range("N1").Value = "DemandMean_50_log1p"
$ws.Range("O1").Value = "DemandMax_50_log1p"
$ws.Range("P1").Value = "DemandMean_25_log1p"
$ws.Range("Q1").Value = "DemandMax_25_log1p"
$ws.Range("R1").Value = "Lag_EoD_Bestand_noSiBe_mean_7Tage"
$ws.Range("S1").Value = "Lag_EoD_Bestand_noSiBe_mean_28Tage"
$ws.Range("T1").Value = "Lag_EoD_Bestand_noSiBe_mean_wbzTage"
$ws.Range("U1").Value = "Lag_EoD_Bestand_noSiBe_mean_2xwbzTage"
$ws.Range("V1").Value = "L_WBZ_BlockMinAbs"

# Fill in G2:V12 with computed feature values
$ws.Range("G2").Value = 2.337
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 3.044522437723423
$ws.Range("J2").Value = 0.09270939810426748
$ws.Range("K2").Value = 2.079441541679836
$ws.Range("L2").Value = 0.097638469563916
$ws.Range("M2").Value = 2.079441541679836
$ws.Range("N2").Value = 0.1076306641923654
$ws.Range("O2").Value = 2.079441541679836
$ws.Range("P2").Value = 0.06595796779179743
$ws.Range("Q2").Value = 1.386294361119891
$ws.Range("R2").Value = -5
$ws.Range("S2").Value = -5
$ws.Range("T2").Value = 3.342857142857143
$ws.Range("U2").Value = 24.08571428571429
$ws.Range("V2").Value = 23

$ws.Range("G3").Value = 2.337
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 3.044522437723423
$ws.Range("J3").Value = 0.167932819341855
$ws.Range("K3").Value = 2.772588722239781
$ws.Range("L3").Value = 0.2076393647782445
$ws.Range("M3").Value = 2.772588722239781
$ws.Range("N3").Value = 0.2500510042341341
$ws.Range("O3").Value = 2.772588722239781
$ws.Range("P3").Value = 0.3429447511268304
$ws.Range("Q3").Value = 2.772588722239781
$ws.Range("R3").Value = -7.142857142857143
$ws.Range("S3").Value = -5.535714285714286
$ws.Range("T3").Value = 3.16
$ws.Range("U3").Value = 23.87142857142857
$ws.Range("V3").Value = 23

$ws.Range("G4").Value = 2.337
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 3.044522437723423
$ws.Range("J4").Value = 0.167932819341855
$ws.Range("K4").Value = 2.772588722239781
$ws.Range("L4").Value = 0.2076393647782445
$ws.Range("M4").Value = 2.772588722239781
$ws.Range("N4").Value = 0.2500510042341341
$ws.Range("O4").Value = 2.772588722239781
$ws.Range("P4").Value = 0.3429447511268304
$ws.Range("Q4").Value = 2.772588722239781
$ws.Range("R4").Value = -9.285714285714286
$ws.Range("S4").Value = -6.071428571428571
$ws.Range("T4").Value = 2.977142857142857
$ws.Range("U4").Value = 23.65714285714286
$ws.Range("V4").Value = 23

$ws.Range("G5").Value = 2.337
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 3.044522437723423
$ws.Range("J5").Value = 0.167932819341855
$ws.Range("K5").Value = 2.772588722239781
$ws.Range("L5").Value = 0.2076393647782445
$ws.Range("M5").Value = 2.772588722239781
$ws.Range("N5").Value = 0.2500510042341341
$ws.Range("O5").Value = 2.772588722239781
$ws.Range("P5").Value = 0.3429447511268304
$ws.Range("Q5").Value = 2.772588722239781
$ws.Range("R5").Value = -11.42857142857143
$ws.Range("S5").Value = -6.607142857142857
$ws.Range("T5").Value = 2.794285714285714
$ws.Range("U5").Value = 23.44285714285714
$ws.Range("V5").Value = 23

$ws.Range("G6").Value = 2.337
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 3.044522437723423
$ws.Range("J6").Value = 0.167932819341855
$ws.Range("K6").Value = 2.772588722239781
$ws.Range("L6").Value = 0.2076393647782445
$ws.Range("M6").Value = 2.772588722239781
$ws.Range("N6").Value = 0.2500510042341341
$ws.Range("O6").Value = 2.772588722239781
$ws.Range("P6").Value = 0.3429447511268304
$ws.Range("Q6").Value = 2.772588722239781
$ws.Range("R6").Value = -13.57142857142857
$ws.Range("S6").Value = -7.142857142857143
$ws.Range("T6").Value = 2.611428571428571
$ws.Range("U6").Value = 23.22857142857143
$ws.Range("V6").Value = 23

$ws.Range("G7").Value = 2.337
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 3.044522437723423
$ws.Range("J7").Value = 0.167932819341855
$ws.Range("K7").Value = 2.772588722239781
$ws.Range("L7").Value = 0.2076393647782445
$ws.Range("M7").Value = 2.772588722239781
$ws.Range("N7").Value = 0.1861022796338607
$ws.Range("O7").Value = 2.772588722239781
$ws.Range("P7").Value = 0.3429447511268304
$ws.Range("Q7").Value = 2.772588722239781
$ws.Range("R7").Value = -15.71428571428571
$ws.Range("S7").Value = -7.678571428571429
$ws.Range("T7").Value = 2.428571428571428
$ws.Range("U7").Value = 23.01428571428572
$ws.Range("V7").Value = 23

$ws.Range("G8").Value = 2.337
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 3.178053830347946
$ws.Range("J8").Value = 0.167932819341855
$ws.Range("K8").Value = 2.772588722239781
$ws.Range("L8").Value = 0.2076393647782445
$ws.Range("M8").Value = 2.772588722239781
$ws.Range("N8").Value = 0.1861022796338607
$ws.Range("O8").Value = 2.772588722239781
$ws.Range("P8").Value = 0.3429447511268304
$ws.Range("Q8").Value = 2.772588722239781
$ws.Range("R8").Value = -17.85714285714286
$ws.Range("S8").Value = -8.214285714285714
$ws.Range("T8").Value = 2.245714285714286
$ws.Range("U8").Value = 22.8
$ws.Range("V8").Value = 23

$ws.Range("G9").Value = 2.337
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 3.178053830347946
$ws.Range("J9").Value = 0.1823215567939546
$ws.Range("K9").Value = 2.772588722239781
$ws.Range("L9").Value = 0.2282586519809802
$ws.Range("M9").Value = 2.772588722239781
$ws.Range("N9").Value = 0.2140110677509373
$ws.Range("O9").Value = 2.772588722239781
$ws.Range("P9").Value = 0.3901976359773759
$ws.Range("Q9").Value = 2.772588722239781
$ws.Range("R9").Value = -20.42857142857143
$ws.Range("S9").Value = -8.857142857142858
$ws.Range("T9").Value = 2.045714285714286
$ws.Range("U9").Value = 22.57714285714286
$ws.Range("V9").Value = 23

$ws.Range("G10").Value = 2.337
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 3.178053830347946
$ws.Range("J10").Value = 0.1823215567939546
$ws.Range("K10").Value = 2.772588722239781
$ws.Range("L10").Value = 0.2282586519809802
$ws.Range("M10").Value = 2.772588722239781
$ws.Range("N10").Value = 0.2140110677509373
$ws.Range("O10").Value = 2.772588722239781
$ws.Range("P10").Value = 0.3901976359773759
$ws.Range("Q10").Value = 2.772588722239781
$ws.Range("R10").Value = -20.85714285714286
$ws.Range("S10").Value = -9.5
$ws.Range("T10").Value = 1.845714285714286
$ws.Range("U10").Value = 22.35428571428572
$ws.Range("V10").Value = 23

$ws.Range("G11").Value = 2.337
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 3.178053830347946
$ws.Range("J11").Value = 0.1823215567939546
$ws.Range("K11").Value = 2.772588722239781
$ws.Range("L11").Value = 0.2282586519809802
$ws.Range("M11").Value = 2.772588722239781
$ws.Range("N11").Value = 0.2140110677509373
$ws.Range("O11").Value = 2.772588722239781
$ws.Range("P11").Value = 0.3901976359773759
$ws.Range("Q11").Value = 2.772588722239781
$ws.Range("R11").Value = -21.28571428571428
$ws.Range("S11").Value = -10.14285714285714
$ws.Range("T11").Value = 1.645714285714286
$ws.Range("U11").Value = 22.13142857142857
$ws.Range("V11").Value = 23

$ws.Range("G12").Value = 2.337
$ws.Range("H12").Value = 4.356708826689592
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0.1823215567939546
$ws.Range("K12").Value = 2.772588722239781
$ws.Range("L12").Value = 0.2282586519809802
$ws.Range("M12").Value = 2.772588722239781
$ws.Range("N12").Value = 0.2140110677509373
$ws.Range("O12").Value = 2.772588722239781
$ws.Range("P12").Value = 0.3901976359773759
$ws.Range("Q12").Value = 2.772588722239781
$ws.Range("R12").Value = -21.71428571428572
$ws.Range("S12").Value = -10.78571428571429
$ws.Range("T12").Value = 1.445714285714286
$ws.Range("U12").Value = 21.90857142857143
$ws.Range("V12").Value = 0

